$d = $word.ActiveDocument

# --- 1) Headline: "Seeking role in asset management firm where I can use my
#        experience in tech..." -> "Seeking PM role where I can use my
#        experience in tech...".  This also removes the old "_GoBack"
#        bookmark that sat between " firm" and " where I can use...".
$rng = $d.Content
$rng.Find.Execute("Seeking role in asset management firm where I can use my experience in tech", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Seeking PM role where I can use my experience in tech", 2)

# --- 2) "...to help the team conduct analysis, build models, and drive
#        profitable investments." -> "...design and scale products,
#        features, and initiatives."
$rng = $d.Content
$rng.Find.Execute("conduct analysis, build models, and drive profitable investments", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "design and scale products, features, and initiatives", 2)

# --- 3) Skills / Business bullet: drop "Marketing, " and "Public Speaking, ",
#        add "Statistics and Analysis, " and swap the tail for "Accounting,
#        Microsoft Excel".
$rng = $d.Content
$rng.Find.Execute("Marketing, Systems Thinking, Lean, Public Speaking, Teambuilding, Strategy, Event Organization", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Systems Thinking, Lean, Teambuilding, Statistics and Analysis, Accounting, Microsoft Excel", 2)

# --- 4) Skills / Finance bullet becomes a Marketing bullet: rename the bold
#        header run and rewrite the list of skills.
$rng = $d.Content
$rng.Find.Execute("Finance", $false, $false, $false, $false, $false, $true, 1, $false, "Marketing", 2)

$rng = $d.Content
$rng.Find.Execute(": Financial Analysis, Accounting, Portfolio Theory, Statistics, Microsoft Excel", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    ": Customer/Developer Relations, Promotion Strategy, Public Speaking, Event Organization", 2)

# --- 5) Re-create the "_GoBack" bookmark inside the new Marketing bullet,
#        matching the authoring tool's placement between "Speaki" and "ng".
$rng = $d.Content
$rng.Find.Execute("Public Speaki")
if ($rng.Find.Found) {
    $bmRange = $d.Range($rng.End, $rng.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
